# usr.xlsx import template: remove the is_locked / is_enabled columns.
#
# Commit intent ("inlineForeignTabs 可编辑聚合表格, isDeleteCascade 级联删除"):
# the is_locked_lbl / is_enabled_lbl generated columns (and their associated
# data-validation select-list template strings) are no longer part of the
# generated header row, so the whole columns are deleted and the columns to
# their right (dept_ids_lbl, role_ids_lbl, rem) shift left to take their
# place, shrinking the header row from 9 columns (A:I) down to 7 (A:G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns E:F hold "<%=comment.is_locked_lbl%>..." and
# "<%=comment.is_enabled_lbl%>..." respectively. Deleting them with a
# shift-left pulls G:I (dept_ids_lbl, role_ids_lbl, rem) into E:G, and the
# now-unreferenced shared strings are dropped on save.
$ws.Range("E1:F1").Delete(-4159)
